# The "Student Details" log sheet keeps a single running row (A2/B2) that
# gets overwritten every time a new login/action event is recorded. This
# commit corresponds to a later run of the program that appended a new
# timestamp entry; the sheet's visible row should now show the latest
# timestamp captured: "02/23/2020 15:32:53".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "02/23/2020 15:32:53"
